# Apply updated cryptos data (prices / % volume changes) to sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "42.341.79"
$ws.Range("E2").Value = "  -0.26%  "
# Row 3
$ws.Range("D3").Value = "2.177.78"
$ws.Range("E3").Value = "  -1.43%  "
# Row 4
$ws.Range("E4").Value = "  +0.02%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "253.74"
$ws.Range("E5").Value = "  +5.53%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.612"
$ws.Range("E6").Value = "  -0.50%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "73.91"
# Row 8
$ws.Range("E8").Value = "  +0.03%  "
# Row 9
$ws.Range("E9").Value = "  -2.68%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.48"
$ws.Range("E10").Value = "  -2.20%  "
# Row 11
$ws.Range("E11").Value = "  -0.98%  "
# Row 12
$ws.Range("E12").Value = "  -0.18%  "
# Row 13
$ws.Range("E13").Value = "  -1.63%  "
# Row 14
$ws.Range("D14").Value = "2.507.45"
$ws.Range("E14").Value = "  -1.14%  "
# Row 15
$ws.Range("E15").Value = "  -3.20%  "
# Row 16
$ws.Range("D16").Value = "2.168.15"
$ws.Range("E16").Value = "  -1.27%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.768"
$ws.Range("E17").Value = "  -3.80%  "
# Row 18
$ws.Range("D18").Value = "42.263.65"
$ws.Range("E18").Value = "  -0.15%  "
# Row 19
$ws.Range("E19").Value = "  -2.54%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.48"
$ws.Range("E20").Value = "  -0.30%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.87"
$ws.Range("E21").Value = "  -0.36%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "226.73"
$ws.Range("E22").Value = "  -0.76%  "
# Row 23
$ws.Range("E23").Value = "  +2.57%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.43"
$ws.Range("E24").Value = "  -7.24%  "
# Row 25
$ws.Range("E25").Value = "  -0.12%  "
# Row 26
$ws.Range("E26").Value = "  -4.40%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.33"
$ws.Range("E27").Value = "  -0.12%  "
# Row 28
$ws.Range("E28").Value = "  +1.61%  "
# Row 29
$ws.Range("E29").Value = "  -2.38%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "170.01"
$ws.Range("E30").Value = "  -1.55%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "36.51"
$ws.Range("E31").Value = "  +9.62%  "
# Row 32
$ws.Range("E32").Value = "  -0.86%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0806"
$ws.Range("E33").Value = "  +2.15%  "
# Row 34
$ws.Range("E34").Value = "  -4.75%  "
# Row 35
$ws.Range("E35").Value = "  -0.91%  "
# Row 36
$ws.Range("E36").Value = "  +0.03%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.24"
$ws.Range("E37").Value = "  -4.03%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0336"
$ws.Range("E38").Value = "  +5.20%  "
# Row 39
$ws.Range("E39").Value = "  -3.39%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.76"
$ws.Range("E40").Value = "  -6.08%  "
# Row 41
$ws.Range("E41").Value = "  +0.54%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "59.31"
$ws.Range("E42").Value = "  -2.21%  "
# Row 43
$ws.Range("E43").Value = "  -6.21%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "102.61"
$ws.Range("E44").Value = "  +3.55%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.468"
$ws.Range("E45").Value = "  +12.80%  "
# Row 46
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.23"
$ws.Range("E46").Value = "  -3.54%  "
# Row 47
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.43"
$ws.Range("E47").Value = "  +6.58%  "
# Row 48
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0967"
$ws.Range("E48").Value = "  -1.22%  "
# Row 49
$ws.Range("E49").Value = "  -1.06%  "
# Row 50
$ws.Range("E50").Value = "  -1.05%  "
# Row 51
$ws.Range("E51").Value = "  +0.43%  "
